$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking value would otherwise be
# auto-converted to a number by Excel (stripping formatting / trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply all cell value updates from the diff, row by row.
$ws.Range("D2").Value = "52.333.40"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.917.42"
$ws.Range("E3").Value = "  +4.79%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "352.03"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").Value = "112.45"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "40.22"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.136"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0861"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("D13").Value = "20.16"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "7.84"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "3.377.65"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "2.920.30"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("D17").Value = "0.995"
$ws.Range("E17").Value = "  +6.31%  "
$ws.Range("D18").Value = "52.413.61"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "3.35"
$ws.Range("E20").Value = "  +5.90%  "
$ws.Range("D21").Value = "14.51"
$ws.Range("E21").Value = "  +8.18%  "
$ws.Range("D22").Value = "0.0₃0982"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "71.15"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "271.58"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "27.04"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "10.62"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").Value = "38.15"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "6.45"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "6.14"
$ws.Range("E33").Value = "  +8.86%  "
$ws.Range("D34").Value = "53.09"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").Value = "0.0935"
$ws.Range("E35").Value = "  +9.36%  "
$ws.Range("D36").Value = "0.0454"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  +6.81%  "
$ws.Range("D39").Value = "18.85"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("D40").Value = "2.08"
$ws.Range("E40").Value = "  +5.57%  "
$ws.Range("E41").Value = "  +9.80%  "
$ws.Range("D42").Value = "24.44"
$ws.Range("E42").Value = "  +12.72%  "
$ws.Range("D44").Value = "122.02"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.214.83"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.56"
$ws.Range("E47").Value = "  +5.08%  "
$ws.Range("E48").Value = "  +6.31%  "
$ws.Range("D49").Value = "0.268"
$ws.Range("E49").Value = "  +25.39%  "
$ws.Range("E50").Value = "  +5.08%  "
$ws.Range("D51").Value = "0.0331"
$ws.Range("E51").Value = "  +14.40%  "
